# Weekly price-sheet update: a new "Ciboulette" price record for
# Femacal de La Calera is inserted as row 137, pushing the existing
# rows 137-211 down to 138-212 (dimension grows from R211 to R212).
#
# Every row in this block shares the same constant columns
# (A, B, C, E, F, G, H, I, K, L, M, N, O, P, Q, R); only the date
# (column D) and the volume (column J) vary per row, so after the
# shift we only need to populate the new row's D/J along with the
# repeated constant values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 137; Excel shifts rows 137:211 down to 138:212
# and grows the used range automatically (A1:R211 -> A1:R212).
$ws.Rows.Item(137).Insert()

$newRow = 137

$ws.Cells.Item($newRow, 1).Value  = 3
$ws.Cells.Item($newRow, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item($newRow, 3).Value  = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value  = 44518
$ws.Cells.Item($newRow, 5).Value  = 5
$ws.Cells.Item($newRow, 6).Value  = 100112039
$ws.Cells.Item($newRow, 7).Value  = "Ciboulette"
$ws.Cells.Item($newRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value  = "Primera"
$ws.Cells.Item($newRow, 10).Value = 160
$ws.Cells.Item($newRow, 11).Value = 1500
$ws.Cells.Item($newRow, 12).Value = 1500
$ws.Cells.Item($newRow, 13).Value = 1500
$ws.Cells.Item($newRow, 14).Value = "`$/docena de atados"
$ws.Cells.Item($newRow, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($newRow, 16).Value = 500
$ws.Cells.Item($newRow, 17).Value = 3
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
